$d = $word.ActiveDocument

# This document's single section has three logo pictures embedded as
# InlineShapes inside the header/footer stories:
#   - Header(2)  -> BTec logo    (wp:docPr/pic:cNvPr name "image2.jpg" -> "image1.jpg")
#   - Footer(1)  -> Pearson logo (wp:docPr/pic:cNvPr name "image1.png" -> "image2.png")
#   - Footer(2)  -> Pearson logo (wp:docPr/pic:cNvPr name "image1.png" -> "image2.png")
#
# The edit simply renames each inline picture (the name shown in Word's
# Selection Pane / exposed as InlineShape.Name), leaving size, alt text,
# id and every other property untouched.

$sec = $d.Sections(1)

$hdrLogo = $sec.Headers(2).Range.InlineShapes(1)
$hdrLogo.Name = "image1.jpg"

$ftr1Logo = $sec.Footers(1).Range.InlineShapes(1)
$ftr1Logo.Name = "image2.png"

$ftr2Logo = $sec.Footers(2).Range.InlineShapes(1)
$ftr2Logo.Name = "image2.png"
